$wb = $excel.ActiveWorkbook
$wsPrev = $wb.Worksheets.Item("Previously added")
$wsNew  = $wb.Worksheets.Item("New")

# ---------------------------------------------------------------------------
# STEP 1: The two rows that were previously in "New" (A2:F3) become the
# newest rows appended at the bottom of "Previously added" (rows 242:243).
# Copy formats first (matching the existing style pattern used throughout
# the sheet), then copy values (this keeps numeric-looking cadastre numbers
# as text, since the source cells are already text/shared-string typed).
# ---------------------------------------------------------------------------

$wsNew.Range("A2:F3").Copy()
$wsPrev.Range("A242:F243").PasteSpecial(-4122) # xlPasteFormats
$wsNew.Range("A2:F3").Copy()
$wsPrev.Range("A242:F243").PasteSpecial(-4163) # xlPasteValues
$excel.CutCopyMode = 0

# Re-create hyperlinks on the moved rows (display text already equals the
# link itself, matching how the rest of the sheet is authored).
$wsPrev.Hyperlinks.Add($wsPrev.Range("A242"), "https://www.ss.com/msg/lv/real-estate/wood/cesis-and-reg/raunas-pag/ndcjg.html", [Type]::Missing, [Type]::Missing, [Type]::Missing)
$wsPrev.Hyperlinks.Add($wsPrev.Range("A243"), "https://www.ss.com/msg/lv/real-estate/wood/limbadzi-and-reg/salacgrivas-l-t/lcnpc.html", [Type]::Missing, [Type]::Missing, [Type]::Missing)

# Hyperlinks.Add re-styles the anchor cell with its own auto "Hyperlink"
# style; restore the original column-A style used throughout the sheet.
$wsPrev.Range("A241").Copy()
$wsPrev.Range("A242").PasteSpecial(-4122) # xlPasteFormats
$wsPrev.Range("A243").PasteSpecial(-4122) # xlPasteFormats
$excel.CutCopyMode = 0

# ---------------------------------------------------------------------------
# STEP 2: Replace the contents of "New" with the four freshly scraped rows.
# ---------------------------------------------------------------------------

# Drop the old hyperlinks/content that used to live in A2:F3.
$wsNew.Range("A2:F3").Hyperlinks.Delete()
$wsNew.Range("A2:F5").ClearContents()

# Give rows 2:5 the same formatting pattern used elsewhere in the workbook
# (style class 3 for the link column, 4 for plain text, 2 for the date).
$wsPrev.Range("A241:F241").Copy()
$wsNew.Range("A2:F2").PasteSpecial(-4122) # xlPasteFormats
$wsNew.Range("A3:F3").PasteSpecial(-4122) # xlPasteFormats
$wsNew.Range("A4:F4").PasteSpecial(-4122) # xlPasteFormats
$wsNew.Range("A5:F5").PasteSpecial(-4122) # xlPasteFormats
$excel.CutCopyMode = 0

$newRows = @(
    @{ Row=2; Link="https://www.ss.com/msg/lv/real-estate/wood/daugavpils-and-reg/naujenes-pag/iokhx.html"; Price="55 000 €"; District="Daugavpils un raj."; Area="4 ha."; Cadastre="44740090231"; Date=45968.709027777775 },
    @{ Row=3; Link="https://www.ss.com/msg/lv/real-estate/wood/kraslava-and-reg/kalniesu-pag/bexeig.html"; Price="10 000 €"; District="Krāslava un raj."; Area="12 ha."; Cadastre="60680040476"; Date=45968.99444444444 },
    @{ Row=4; Link="https://www.ss.com/msg/lv/real-estate/wood/madona-and-reg/metrienas-pag/pglli.html"; Price="5 000 €"; District="Madona un raj."; Area="2 ha."; Cadastre="70760050003"; Date=45970.81180555555 },
    @{ Row=5; Link="https://www.ss.com/msg/lv/real-estate/wood/rezekne-and-reg/cornajas-pag/ikbfo.html"; Price="68 000 €"; District="Rēzekne un raj."; Area="11.50 ha."; Cadastre="78460090131"; Date=45969.950694444444 }
)

foreach ($r in $newRows) {
    $row = $r.Row
    # The cadastre numbers are purely numeric strings; force the cell to
    # text format first so Excel doesn't silently store them as numbers.
    $wsNew.Range("E$row").NumberFormat = "@"
    $wsNew.Range("A$row").Value = $r.Link
    $wsNew.Range("B$row").Value = $r.Price
    $wsNew.Range("C$row").Value = $r.District
    $wsNew.Range("D$row").Value = $r.Area
    $wsNew.Range("E$row").Value = $r.Cadastre
    $wsNew.Range("F$row").Value = $r.Date
}

# Re-apply the shared formatting once more: assigning NumberFormat above
# nudged column E onto its own style, and Hyperlinks.Add (next) will nudge
# column A onto its own style too, so do the final format pass afterwards.
foreach ($r in $newRows) {
    $row = $r.Row
    $wsNew.Hyperlinks.Add($wsNew.Range("A$row"), $r.Link, [Type]::Missing, [Type]::Missing, [Type]::Missing)
}

$wsPrev.Range("A241:F241").Copy()
$wsNew.Range("A2:F2").PasteSpecial(-4122) # xlPasteFormats
$wsNew.Range("A3:F3").PasteSpecial(-4122) # xlPasteFormats
$wsNew.Range("A4:F4").PasteSpecial(-4122) # xlPasteFormats
$wsNew.Range("A5:F5").PasteSpecial(-4122) # xlPasteFormats
$excel.CutCopyMode = 0

Write-Host "Edit applied"
